$d = $word.ActiveDocument

# Locate the full span covering the three original runs:
#   " requirements of Paragraph 45 of the HUD agreement. The consolidation reported insufficient staff to meet th"
# + "is portion "
# + "of the requirements. "
$old = " requirements of Paragraph 45 of the HUD agreement. The consolidation reported insufficient staff to meet this portion of the requirements. "
$full = $d.Content.Text
$start = $full.IndexOf($old)
$end = $start + $old.Length

$span = $d.Range($start, $end)

# Replace with the new combined text: truncated sentence + trailing single space.
$newText = " requirements of Paragraph 45 of the HUD agreement. "
$span.Text = $newText

function Toggle-Split($r) {
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# Force the edited text to become its own run, detached from the preceding
# "...and collection" run (which must stay untouched).
$editedRange = $d.Range($start, $start + $newText.Length)
Toggle-Split $editedRange

# Now split off the trailing single space into its own run, detached from
# the truncated-sentence run.
$spaceStart = $start + ($newText.Length - 1)
$spaceEnd = $start + $newText.Length
$spaceRange = $d.Range($spaceStart, $spaceEnd)
Toggle-Split $spaceRange
